$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.864.00'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').Value = '2.457.14'
$ws.Range('E3').Value = '  -2.39%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '578.80'
$ws.Range('E5').Value = '  -2.61%  '
$ws.Range('D6').Value = '165.65'
$ws.Range('E6').Value = '  -4.94%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -3.05%  '
$ws.Range('D9').Value = '2.453.68'
$ws.Range('E9').Value = '  -2.54%  '
$ws.Range('E10').Value = '  -4.37%  '
$ws.Range('E11').Value = '  -1.13%  '
$ws.Range('D12').Value = '4.87'
$ws.Range('E12').Value = '  -4.57%  '
$ws.Range('D13').Value = '0.332'
$ws.Range('E13').Value = '  -3.45%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.900.18'
$ws.Range('E14').Value = '  -2.56%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = '25.30'
$ws.Range('E15').Value = '  -4.47%  '
$ws.Range('D16').Value = '66.808.24'
$ws.Range('E16').Value = '  -1.48%  '
$ws.Range('E17').Value = '  -5.60%  '
$ws.Range('D18').Value = '2.442.63'
$ws.Range('E18').Value = '  -1.70%  '
$ws.Range('D19').Value = '11.33'
$ws.Range('E19').Value = '  -4.66%  '
$ws.Range('D20').Value = '7.71'
$ws.Range('E20').Value = '  -3.43%  '
$ws.Range('D21').Value = '354.56'
$ws.Range('E21').Value = '  -2.60%  '
$ws.Range('D22').Value = '4.01'
$ws.Range('E22').Value = '  -3.29%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').Value = '69.35'
$ws.Range('E24').Value = '  -2.80%  '
$ws.Range('D25').Value = '4.21'
$ws.Range('E25').Value = '  -9.07%  '
$ws.Range('E26').Value = '  -8.08%  '
$ws.Range('D27').Value = '8.89'
$ws.Range('E27').Value = '  -10.77%  '
$ws.Range('D28').Value = '0.996'
$ws.Range('E28').Value = '  -0.44%  '
$ws.Range('D29').Value = '2.576.41'
$ws.Range('E29').Value = '  -2.49%  '
$ws.Range('D30').Value = '0.0₃0898'
$ws.Range('E30').Value = '  -8.32%  '
$ws.Range('D31').Value = '505.41'
$ws.Range('E31').Value = '  -4.86%  '
$ws.Range('D32').Value = '7.78'
$ws.Range('E32').Value = '  -6.45%  '
$ws.Range('E33').Value = '  -6.22%  '
$ws.Range('E34').Value = '  -8.23%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').Value = '159.05'
$ws.Range('E36').Value = '  +0.56%  '
$ws.Range('E37').Value = '  -10.11%  '
$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').Value = '18.45'
$ws.Range('E38').Value = '  -1.44%  '
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').Value = '18.56'
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('E40').Value = '  -7.29%  '
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('D42').Value = '1.66'
$ws.Range('E42').Value = '  -7.29%  '
$ws.Range('E43').Value = '  -7.14%  '
$ws.Range('D44').Value = '4.76'
$ws.Range('E44').Value = '  -7.72%  '
$ws.Range('D45').Value = '38.68'
$ws.Range('E45').Value = '  -3.12%  '
$ws.Range('E46').Value = '  -7.88%  '
$ws.Range('D47').Value = '141.25'
$ws.Range('E47').Value = '  -4.04%  '
$ws.Range('E48').Value = '  -6.29%  '
$ws.Range('D49').Value = '0.514'
$ws.Range('E49').Value = '  -7.11%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0251'
$ws.Range('E50').Value = '  -9.27%  '
$ws.Range('B51').Value = 'Optimism'
$ws.Range('C51').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D51').Value = '1.58'
$ws.Range('E51').Value = '  -8.30%  '
